$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{Row=2; C=38485; D=55646015},
    @{Row=3; C=92369; D=135399008},
    @{Row=4; C=31588; D=46778372},
    @{Row=5; C=8854; D=13159856},
    @{Row=6; C=2044; D=3037971},
    @{Row=7; C=164; D=241093},
    @{Row=12; C=41935; D=56879024},
    @{Row=13; C=9829; D=14217977},
    @{Row=14; C=26317; D=38588225},
    @{Row=15; C=8407; D=12475978},
    @{Row=16; C=2189; D=3253430},
    @{Row=20; C=10357; D=13701648},
    @{Row=21; C=13599; D=19625866},
    @{Row=22; C=32061; D=47047943},
    @{Row=23; C=10339; D=15368810},
    @{Row=24; C=2674; D=3975771},
    @{Row=27; C=11857; D=15830875},
    @{Row=28; C=7813; D=11310704},
    @{Row=29; C=22859; D=33552390},
    @{Row=30; C=7905; D=11757591},
    @{Row=31; C=1999; D=2982919},
    @{Row=34; C=8433; D=11140577},
    @{Row=35; C=3325; D=4802652},
    @{Row=36; C=7979; D=11652940},
    @{Row=37; C=3219; D=4771461},
    @{Row=41; C=2517; D=3401422},
    @{Row=42; C=17565; D=25399136},
    @{Row=43; C=51900; D=76075091},
    @{Row=44; C=19213; D=28535892},
    @{Row=45; C=5695; D=8477760},
    @{Row=46; C=1231; D=1837045},
    @{Row=50; C=17008; D=22611560},
    @{Row=51; C=2105; D=3053378},
    @{Row=52; C=7130; D=10477142},
    @{Row=53; C=2409; D=3597964},
    @{Row=54; C=769; D=1148805},
    @{Row=55; C=194; D=287226},
    @{Row=57; C=7230; D=9939333},
    @{Row=58; C=1102; D=1809629},
    @{Row=59; C=2730; D=4487158},
    @{Row=60; C=1074; D=1766338},
    @{Row=61; C=367; D=606883},
    @{Row=62; C=121; D=205100},
    @{Row=63; C=23; D=42000},
    @{Row=64; C=1615; D=2468762},
    @{Row=65; C=15666; D=22625567},
    @{Row=66; C=45386; D=66407618},
    @{Row=67; C=15905; D=23631320},
    @{Row=68; C=4625; D=6888551},
    @{Row=69; C=951; D=1414668},
    @{Row=73; C=15305; D=20160900},
    @{Row=74; C=53485; D=77835879},
    @{Row=75; C=150479; D=221691137},
    @{Row=76; C=65049; D=96930306},
    @{Row=77; C=20813; D=31099822},
    @{Row=78; C=4963; D=7413401},
    @{Row=85; C=52765; D=71721377},
    @{Row=86; C=4725; D=6847767},
    @{Row=87; C=11793; D=17323302},
    @{Row=88; C=3938; D=5869583},
    @{Row=89; C=1363; D=2036989},
    @{Row=90; C=290; D=432512},
    @{Row=93; C=5515; D=7413715},
    @{Row=94; C=1636; D=2356862},
    @{Row=95; C=5292; D=7795801},
    @{Row=96; C=1968; D=2930326},
    @{Row=97; C=702; D=1051960},
    @{Row=101; C=3645; D=4823155},
    @{Row=102; C=690; D=1125325},
    @{Row=103; C=416; D=695597},
    @{Row=104; C=153; D=254020},
    @{Row=107; C=10995; D=15951287},
    @{Row=108; C=29602; D=43481409},
    @{Row=109; C=9905; D=14726705},
    @{Row=110; C=2726; D=4064080},
    @{Row=114; C=9937; D=13122964},
    @{Row=115; C=31014; D=44719206},
    @{Row=116; C=67068; D=98140593},
    @{Row=117; C=21620; D=32131108},
    @{Row=118; C=6131; D=9134521},
    @{Row=119; C=1147; D=1714100},
    @{Row=124; C=26208; D=34980794},
    @{Row=125; C=36715; D=52980911},
    @{Row=126; C=78004; D=114053645},
    @{Row=127; C=24148; D=35841592},
    @{Row=128; C=6481; D=9631358},
    @{Row=129; C=1268; D=1885311},
    @{Row=133; C=32278; D=42842162},
    @{Row=134; C=13520; D=19570946},
    @{Row=135; C=32782; D=48143190},
    @{Row=136; C=11617; D=17259587},
    @{Row=137; C=3005; D=4478741},
    @{Row=138; C=508; D=755990},
    @{Row=141; C=10952; D=14600633},
    @{Row=142; C=35825; D=51741254},
    @{Row=143; C=82691; D=121145013},
    @{Row=144; C=24709; D=36709056},
    @{Row=145; C=6485; D=9677067},
    @{Row=146; C=1469; D=2185730},
    @{Row=149; C=29667; D=39996156}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
}
